$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab name acts as the file name label here)
$ws.Name = "前端keyValue1696814218288.xlsx"

# Add the two new rows of key/value translations
$ws.Range("A2").Value = "ce5b42"
$ws.Range("B2").Value = "天才/帅哥。。"

$ws.Range("A3").Value = "ce5b42-Loading"
$ws.Range("B3").Value = "加载中..."

# Keep the "numbers stored as text" error-check ignore flag covering the
# whole used range (A1:C1 -> A1:C3) now that more rows were added.
$ws.Range("A1:C3").Errors.Item(1).Ignore = $true
